$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.831.25"
$ws.Range("E2").Value = "  +1.38%  "

$ws.Range("D3").Value = "2.634.82"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.38%  "

$ws.Range("E10").Value = "  +1.86%  "

$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("E12").Value = "  -1.64%  "

$ws.Range("D13").Value = "3.104.51"
$ws.Range("E13").Value = "  +2.06%  "

$ws.Range("D14").Value = "59.743.88"
$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.59%  "

$ws.Range("D16").Value = "2.667.44"
$ws.Range("E16").Value = "  +3.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "343.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.11%  "

$ws.Range("E24").Value = "  +2.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  +3.19%  "

$ws.Range("D28").Value = "0.0₃0753"
$ws.Range("E28").Value = "  +5.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("E34").Value = "  +1.81%  "

$ws.Range("E35").Value = "  +1.95%  "

$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.840"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.14%  "

$ws.Range("E38").Value = "  +2.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "291.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0533"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.94%  "

$ws.Range("D46").Value = "1.967.79"
$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("E47").Value = "  +1.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.51%  "

$ws.Range("E49").Value = "  +2.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.01%  "
